$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H33").Value = 585.35297
$ws.Range("I33").Value = 144.28572
$ws.Range("J33").Value = 894.1
$ws.Range("K33").Value = 144.28572
$ws.Range("L33").Value = 894.1
$ws.Range("M33").Value = 84.71428
$ws.Range("N33").Value = -1352.1
$ws.Range("H41").Value = 1151.3334
$ws.Range("I41").Value = 288.8
$ws.Range("K41").Value = 288.8
$ws.Range("M41").Value = 151.2
$ws.Range("H43").Value = 2614.5
$ws.Range("J43").Value = 2424.5
$ws.Range("L43").Value = 2424.5
$ws.Range("N43").Value = -2562.5
$ws.Range("H55").Value = 352.4
$ws.Range("J55").Value = 999.5
$ws.Range("L55").Value = 999.5
$ws.Range("N55").Value = -1427.5
$ws.Range("H74").Value = 2899
$ws.Range("I74").Value = 0
$ws.Range("K74").Value = 0
$ws.Range("M74").ClearContents()
$ws.Range("H77").Value = 2899
$ws.Range("I77").Value = 0
$ws.Range("K77").Value = 0
$ws.Range("M77").ClearContents()
$ws.Range("H101").Value = 761.4286
$ws.Range("I101").Value = 729
$ws.Range("J101").Value = 842.5
$ws.Range("K101").Value = 2187
$ws.Range("L101").Value = 2527.5
$ws.Range("M101").Value = -565
$ws.Range("N101").Value = -5771.5
$ws.Range("H116").Value = 12452.363
$ws.Range("I116").Value = 13886.556
$ws.Range("K116").Value = 13886.556
$ws.Range("M116").Value = -10444.556
$ws.Range("H127").Value = 5415.2856
$ws.Range("I127").Value = 1546.2858
$ws.Range("K127").Value = 4638.857400000001
$ws.Range("M127").Value = 321.1425999999992
$ws.Range("H137").Value = 3747343.8
$ws.Range("I137").Value = 5949440
$ws.Range("K137").Value = 17848320
$ws.Range("M137").Value = -17845770
$ws.Range("H139").Value = 0
$ws.Range("J139").Value = 0
$ws.Range("L139").Value = 0
$ws.Range("N139").ClearContents()
$ws.Range("H141").Value = 5766
$ws.Range("J141").Value = 8266
$ws.Range("L141").Value = 24798
$ws.Range("N141").Value = -35158
$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H102").Value = 2150.8462
$ws.Range("I102").Value = 1851.381
$ws.Range("K102").Value = 1851.381
$ws.Range("M102").Value = -229.3810000000001
$ws.Range("H124").Value = 48714.5
$ws.Range("J124").Value = 48714.5
$ws.Range("L124").Value = 48714.5
$ws.Range("N124").Value = -58534.5
$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H20").Value = 1552.7097
$ws.Range("I20").Value = 1553.8182
$ws.Range("K20").Value = 1553.8182
$ws.Range("M20").Value = -1306.8182
$ws.Range("H105").Value = 2031.8235
$ws.Range("J105").Value = 1759
$ws.Range("L105").Value = 1759
$ws.Range("N105").Value = -5253
$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H4").Value = 8490
$ws.Range("J4").Value = 8490
$ws.Range("L4").Value = 8490
$ws.Range("N4").Value = -8714
$ws.Range("H134").Value = 3830.628
$ws.Range("I134").Value = 4174.108
$ws.Range("K134").Value = 12522.324
$ws.Range("M134").Value = -9987.324000000001
$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H11").Value = 1006.5
$ws.Range("I11").Value = 1504.75
$ws.Range("K11").Value = 4514.25
$ws.Range("M11").Value = -4374.25
$ws.Range("H44").Value = 2236.24
$ws.Range("I44").Value = 917.1667
$ws.Range("J44").Value = 2652.7896
$ws.Range("K44").Value = 2751.5001
$ws.Range("L44").Value = 7958.3688
$ws.Range("M44").Value = -2353.5001
$ws.Range("N44").Value = -8754.3688
$ws.Range("H55").Value = 6080.143
$ws.Range("J55").Value = 9834.5
$ws.Range("L55").Value = 29503.5
$ws.Range("N55").Value = -29857.5
$ws.Range("H68").Value = 16671509
$ws.Range("J68").Value = 7902.5
$ws.Range("L68").Value = 23707.5
$ws.Range("N68").Value = -25329.5
$ws.Range("H71").Value = 16671509
$ws.Range("J71").Value = 7902.5
$ws.Range("L71").Value = 71122.5
$ws.Range("N71").Value = -79234.5
$ws.Range("H131").Value = 5377787
$ws.Range("I131").Value = 166667600
$ws.Range("J131").Value = 1459.8556
$ws.Range("K131").Value = 500002800
$ws.Range("L131").Value = 4379.566800000001
$ws.Range("M131").Value = -499997760
$ws.Range("N131").Value = -14459.5668
$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H70").Value = 8263.723
$ws.Range("J70").Value = 6625.7144
$ws.Range("L70").Value = 6625.7144
$ws.Range("N70").Value = -7165.7144
$ws.Range("H73").Value = 8263.723
$ws.Range("J73").Value = 6625.7144
$ws.Range("L73").Value = 6625.7144
$ws.Range("N73").Value = -8497.714400000001
$ws.Range("H122").Value = 2312.6155
$ws.Range("I122").Value = 2343.6667
$ws.Range("K122").Value = 7031.000100000001
$ws.Range("M122").Value = -4581.000100000001
$ws.Range("H126").Value = 6211.7856
$ws.Range("I126").Value = 9012.166999999999
$ws.Range("J126").Value = 4111.5
$ws.Range("K126").Value = 27036.501
$ws.Range("L126").Value = 12334.5
$ws.Range("M126").Value = -24566.501
$ws.Range("N126").Value = -17274.5
$ws.Range("H132").Value = 40277.223
$ws.Range("I132").Value = 43554.816
$ws.Range("K132").Value = 130664.448
$ws.Range("M132").Value = -128134.448
$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H16").Value = 839.6
$ws.Range("I16").Value = 898.6667
$ws.Range("K16").Value = 898.6667
$ws.Range("M16").Value = -728.6667
$ws.Range("H40").Value = 6849.478
$ws.Range("I40").Value = 6519.5
$ws.Range("K40").Value = 6519.5
$ws.Range("M40").Value = -6383.5
$ws.Range("H93").Value = 83334990
$ws.Range("I93").Value = 938.2222
$ws.Range("J93").Value = 333337150
$ws.Range("K93").Value = 938.2222
$ws.Range("L93").Value = 333337150
$ws.Range("M93").Value = 309.7778
$ws.Range("N93").Value = -333339646
$ws.Range("H122").Value = 6084.35
$ws.Range("I122").Value = 5965.6665
$ws.Range("K122").Value = 17896.9995
$ws.Range("M122").Value = -15446.9995
$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H5").Value = 550
$ws.Range("J5").Value = 550
$ws.Range("L5").Value = 550
$ws.Range("N5").Value = -774
$ws.Range("H23").Value = 3750
$ws.Range("J23").Value = 4500
$ws.Range("L23").Value = 4500
$ws.Range("N23").Value = -4958
$ws.Range("H54").Value = 24995
$ws.Range("J54").Value = 24995
$ws.Range("L54").Value = 24995
$ws.Range("N54").Value = -26035
$ws.Range("H113").Value = 1254.3334
$ws.Range("I113").Value = 1338.1
$ws.Range("J113").Value = 1149.625
$ws.Range("K113").Value = 4014.3
$ws.Range("L113").Value = 3448.875
$ws.Range("M113").Value = -1844.3
$ws.Range("N113").Value = -7788.875
$ws.Range("H122").Value = 2043.0555
$ws.Range("I122").Value = 1991.8
$ws.Range("K122").Value = 5975.4
$ws.Range("M122").Value = -3525.4
$ws.Range("H132").Value = 2880.1052
$ws.Range("I132").Value = 6111
$ws.Range("K132").Value = 18333
$ws.Range("M132").Value = -15803
